$wb = $excel.ActiveWorkbook

# Update SimParameters sheet values
$sim = $wb.Worksheets.Item("SimParameters")

# Swap treatment-effect multipliers (B4 <-> B6)
$sim.Range("B4").Value = 4
$sim.Range("B6").Value = 1.33

# Swap abortion severity multipliers (B8 <-> B10)
$sim.Range("B8").Value = 0.25
$sim.Range("B10").Value = 0.75

# Make SimParameters the active/selected sheet with B11 selected
$sim.Activate()
$sim.Range("B11").Select()
